$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cesar")

$values = @{
    3  = @(10, 20)
    4  = @(5, 35)
    5  = @(20, 35)
    6  = @(65, 10)
    9  = @(20, 30)
    10 = @(10, 30)
    11 = @(5, 10)
    12 = @(65, 30)
    15 = @(20, 25)
    16 = @(5, 25)
    17 = @(10, 25)
    18 = @(65, 25)
    21 = @(5, 30)
    22 = @(10, 30)
    23 = @(20, 30)
    24 = @(65, 10)
    27 = @(50, 30)
    28 = @(10, 20)
    29 = @(10, 25)
    30 = @(30, 25)
    33 = @(30, 30)
    34 = @(10, 30)
    35 = @(20, 10)
    36 = @(40, 30)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Range("D$row").Value = $pair[0]
    $ws.Range("E$row").Value = $pair[1]
}

$consolidated = $wb.Worksheets.Item("Consolidated")
$consolidated.Range("H10").Formula = "=H9/24"
$consolidated.Range("I10").Formula = "=I9/24"
$consolidated.Range("J10").Formula = "=J9/24"
$consolidated.Range("K10").Formula = "=K9/24"
$consolidated.Range("N10").Formula = "=N9/24"
$consolidated.Range("O10").Formula = "=O9/24"
$consolidated.Range("P10").Formula = "=P9/24"
$consolidated.Range("Q10").Formula = "=Q9/24"

$carley = $wb.Worksheets.Item("Carley")
[void]$carley.Activate()
[void]$carley.Range("D3").Select()

[void]$consolidated.Activate()
[void]$consolidated.Range("I13").Select()

[void]$ws.Activate()
[void]$ws.Range("D3").Select()

Write-Host "values set"
